$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 38, shifting existing rows 38:60 down to 39:61
$ws.Rows("38").Insert()

# Populate the newly inserted row 38 with the new record
$ws.Range("A38").Value = 4
$ws.Range("B38").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C38").Value = "Los Lagos"
$ws.Range("D38").Value = (Get-Date -Year 2023 -Month 2 -Day 3).Date
$ws.Range("E38").Value = 10
$ws.Range("F38").Value = 100112030
$ws.Range("G38").Value = "Poroto granado"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 80
$ws.Range("K38").Value = 45000
$ws.Range("L38").Value = 45000
$ws.Range("M38").Value = 45000
$ws.Range("N38").Value = "`$/saco 25 kilos"
$ws.Range("O38").Value = "Región Metropolitana"
$ws.Range("P38").Value = 1800
$ws.Range("Q38").Value = 25
$ws.Range("R38").Value = "Hortaliza"
